$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_2_06_B")

# ---------------------------------------------------------------
# 2. Update the report title / period text (shared strings)
# ---------------------------------------------------------------
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Thousand Tons)"

# ---------------------------------------------------------------
# 3. Insert a new data row for "November" (2016) just above the
#    "Year to Date" summary block (old row 53 -> new row 54+)
# ---------------------------------------------------------------
$ws.Rows.Item(53).Insert()
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 161
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 51
$ws.Range("E53").Value = 96
$ws.Range("F53").Value = 14

# ---------------------------------------------------------------
# 4. Refresh the "Year to Date" annual totals (now rows 55-57)
# ---------------------------------------------------------------
$ws.Range("B55").Value = 1786
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 591
$ws.Range("E55").Value = 1011
$ws.Range("F55").Value = 184

$ws.Range("B56").Value = 1812
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 595
$ws.Range("E56").Value = 1030
$ws.Range("F56").Value = 186

$ws.Range("B57").Value = 1805
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 606
$ws.Range("E57").Value = 1046
$ws.Range("F57").Value = 152

# ---------------------------------------------------------------
# 5. Update the "Rolling 12 Months Ending in ..." header text and
#    its totals (now rows 58-60)
# ---------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$ws.Range("B59").Value = 1981
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 654
$ws.Range("E59").Value = 1124
$ws.Range("F59").Value = 203

$ws.Range("B60").Value = 1979
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 667
$ws.Range("E60").Value = 1143
$ws.Range("F60").Value = 169
